$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 706, pushing existing rows 706-747 down to 707-748.
$ws.Rows(706).Insert()

# Write the new data point: 2026/01/22 (Thursday), hour 23, ranking 161.
# Force column A to be treated as text so the date string is stored verbatim
# (matching the rest of the sheet, which stores dates as plain text), then
# reset the style back to the sheet's default so no extra formatting sticks.
$ws.Range("A706").NumberFormat = "@"
$ws.Range("A706").Value = "2026/01/22"
$ws.Range("A706").Style = "Normal"

$ws.Range("B706").Value = "木"
$ws.Range("C706").Value = 23
$ws.Range("D706").Value = 161
